# 4.c.1 — add a new "2023" data column (Q) to the right of the existing
# "2022" column (P), carrying the same per-row formatting, then reset the
# sheet's active cell/selection back to the default (A1) instead of the
# stale "Q4" selection left over from editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (thin divider row above the year header, no values) ---
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 4 (year header row) : Q4 = 2023 ---
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2023

# --- Row 5 (a) preschool) : Q5 = 93.7 ---
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 93.7

# --- Row 6 (b) primary school) : Q6 = 95.5 ---
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 95.5

# --- Row 7 (c/d main & high school) : Q7 = 97.1 ---
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 97.1

$excel.CutCopyMode = $false

# Reset the selection that was parked on Q4 back to the sheet's default
# top-left cell.
$ws.Range("A1").Select()
